$wb = $excel.ActiveWorkbook

# Helper to write a value as TEXT even when it looks numeric, mirroring how
# these generator output cells are stored as text in the workbook.
function Set-TextValue($ws, $cellRef, [string]$text) {
    $rng = $ws.Range($cellRef)
    if ($text -match '^[\+\-]?(\d+\.?\d*|\.\d+)([eE][\+\-]?\d+)?$') {
        # Purely numeric-looking string: prefix with an apostrophe so Excel
        # keeps it stored as text instead of converting it to a number.
        $rng.Value = "'" + $text
    } else {
        $rng.Value = $text
    }
}

# Restricciones_del_follower sheet (3rd sheet)
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws "A2" "5.95 - y"
Set-TextValue $ws "B2" "-5.95"
Set-TextValue $ws "D2" "0.37"
Set-TextValue $ws "E2" "7.9"
Set-TextValue $ws "F2" "4.4"

Set-TextValue $ws "A3" "-0.6000000000000005 - x + y"
Set-TextValue $ws "B3" "-2.3999999999999995"
Set-TextValue $ws "D3" "0.44"
Set-TextValue $ws "E3" "3.5999999999999996"
Set-TextValue $ws "F3" "5.6000000000000005"

Set-TextValue $ws "A4" "-17.25 + x + 2y"
Set-TextValue $ws "B4" "5.25"
Set-TextValue $ws "D4" "0.0"
Set-TextValue $ws "E4" "2.6"
Set-TextValue $ws "F4" "5.300000000000001"

Set-TextValue $ws "A5" "-15.69 + 4x - y"
Set-TextValue $ws "B5" "3.4499999999999993"
Set-TextValue $ws "D5" "0.07"
Set-TextValue $ws "E5" "4.699999999999999"
Set-TextValue $ws "F5" "9.399999999999999"

# Punto_modificado sheet (4th sheet)
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws "A2" "5.35"
Set-TextValue $ws "B2" "5.95"

# Vector_bf sheet (5th sheet)
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "-1.0"

# Vector_BF sheet (6th sheet)
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "-16.799999999999997"
Set-TextValue $ws "A3" "6.8"
